$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data
#    rows (2..304) from 2023-09-23 (45192) to 2023-10-03 (45202).
$ws.Range("C2:C304").Value = 45202

# 2. Row 304 gains an explicit row height (15pt, customHeight) in the
#    canonical XML — setting RowHeight explicitly reproduces that.
$ws.Rows.Item(304).RowHeight = 15

# 3. Append the new record as row 305.
$r = 305
$ws.Cells.Item($r, 1).Value = "A 45845-2023"

$ws.Cells.Item($r, 2).Value = 45195
$ws.Cells.Item($r, 2).NumberFormat = $ws.Cells.Item(304, 2).NumberFormat()

$ws.Cells.Item($r, 3).Value = 45202
$ws.Cells.Item($r, 3).NumberFormat = $ws.Cells.Item(304, 3).NumberFormat()

$ws.Cells.Item($r, 4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item($r, 5).Value = "ÖVERKALIX"
$ws.Cells.Item($r, 6).Value = "Sveaskog"

$ws.Cells.Item($r, 7).Value = 4
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0
$ws.Cells.Item($r, 11).Value = 0
$ws.Cells.Item($r, 12).Value = 0
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 0
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 0

# Column R keeps the wrap-text style used throughout the sheet, with no
# content (empty inline string cell).
$ws.Cells.Item($r, 18).Value = ""
$ws.Cells.Item($r, 18).WrapText = $true
